# add error bars to observed rates in plots
#
# Functional changes captured from the source diff:
#   1. Column T (rows 18-28, "error bar" / observed-rate deviation column)
#      is re-entered as a single fill-down formula series, turning it into
#      a shared formula group anchored at T18 (T18:T28) - mirrors the
#      author dragging the T18 formula down to T28.
#   2. Several "max" recruitment counts in column M (rows 22-28) were
#      corrected, which ripples into the recruitment-rate formulas in
#      column N (AVERAGE(L:M)/AVERAGE(J:K)).
#   3. The active selection moved from V34 to O19.
#   4. The workbook window was shifted to a second monitor (xWindow).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MinCount_summary_KZ-withimm")

# --- 1. Re-create column T (T18:T28) as one fill-down formula series ---
# T18 is the "source" cell; T19:T28 are filled down from it (relative refs
# shift per-row, exactly like Excel's own fill-down / shared formula).
for ($r = 18; $r -le 28; $r++) {
    $ws.Cells.Item($r, 20).Formula = "=(AVERAGE(J${r}:K${r})/(Q${r}-R${r}))-S${r}"
}

# --- 2. Corrected "max" values in column M (rows 22-28) ---
$ws.Range("M22").Value = 25
$ws.Range("M23").Value = 29
$ws.Range("M24").Value = 31
$ws.Range("M25").Value = 27
$ws.Range("M26").Value = 39
$ws.Range("M27").Value = 30
$ws.Range("M28").Value = 38

# --- 3. Move the active selection to O19 ---
$ws.Activate() | Out-Null
$ws.Range("O19").Select() | Out-Null

# --- 4. Shift workbook window position (best effort; matches xWindow) ---
$wb.Windows.Item(1).Left = -40800

$wb.Save() | Out-Null
